$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.936.15"
$ws.Range("E2").Value = "  +1.74%  "

# Row 3
$ws.Range("D3").Value = "3.334.84"
$ws.Range("E3").Value = "  +1.80%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.78"
$ws.Range("E5").Value = "  +1.92%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.50"
$ws.Range("E6").Value = "  +0.94%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +2.00%  "

# Row 9
$ws.Range("D9").Value = "3.330.63"
$ws.Range("E9").Value = "  +1.80%  "

# Row 10
$ws.Range("E10").Value = "  +4.67%  "

# Row 11
$ws.Range("E11").Value = "  +1.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.84"
$ws.Range("E12").Value = "  +2.70%  "

# Row 13
$ws.Range("E13").Value = "  +1.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "689.02"
$ws.Range("E14").Value = "  -1.26%  "

# Row 15
$ws.Range("D15").Value = "3.872.09"
$ws.Range("E15").Value = "  +1.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.44"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17
$ws.Range("D17").Value = "67.948.22"
$ws.Range("E17").Value = "  +1.58%  "

# Row 18
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").Value = "3.331.03"
$ws.Range("E19").Value = "  +1.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.41"
$ws.Range("E20").Value = "  +0.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.07"
$ws.Range("E21").Value = "  +3.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.897"
$ws.Range("E22").Value = "  +1.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.05"
$ws.Range("E23").Value = "  +0.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.33"
$ws.Range("E24").Value = "  +3.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.22"
$ws.Range("E25").Value = "  +1.05%  "

# Row 26
$ws.Range("E26").Value = "  +0.39%  "

# Row 27
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  +2.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.97"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("E30").Value = "  +1.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("E31").Value = "  +4.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "574.55"
$ws.Range("E32").Value = "  +1.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.01"
$ws.Range("E33").Value = "  +1.88%  "

# Row 34
$ws.Range("E34").Value = "  +1.96%  "

# Row 35
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.708.97"
$ws.Range("E36").Value = "  -4.69%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.87"
$ws.Range("E37").Value = "  +2.49%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  -2.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.55"
$ws.Range("E39").Value = "  +8.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  +2.68%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").Value = "  +2.48%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.21"
$ws.Range("E42").Value = "  +6.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.37"
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$ws.Range("D44").Value = "0.0₃0678"
$ws.Range("E44").Value = "  +0.64%  "

# Row 45
$ws.Range("E45").Value = "  +2.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("E46").Value = "  +0.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  +6.10%  "

# Row 48
$ws.Range("E48").Value = "  +1.14%  "

# Row 49
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.32"
$ws.Range("E50").Value = "  -2.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.75"
$ws.Range("E51").Value = "  +0.03%  "
